# PNAD 2009 - roubo (cv124153a)
# Correção nos dados: a linha 6 ("grandes regiões e unidades da federação")
# era um cabeçalho de seção sem valores numéricos e foi removida por engano
# durante a geração da planilha; os dados de cada UF estavam, na verdade,
# deslocados uma linha abaixo de onde deveriam estar. Removendo a linha 6
# inteira faz com que todas as linhas subsequentes (Norte, Rondônia, Acre,
# ..., Goiás) subam uma posição, corrigindo o alinhamento entre os rótulos
# da coluna A e os valores numéricos das colunas B:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
